{"js": "// Replace the date line and each of the division problems in the table\n// with their updated values, as described by the diff. Each lookup\n// string is unique in the document, so a plain search + replace for\n// each pair is safe and unambiguous (even though some replacement\n// values repeat, e.g. \"55\u00f78=\" appears twice).\nconst replacements = [\n  [\"2025-08-05 Tuesday\", \"2025-08-06 Wednesday\"],\n  [\"40\u00f74=\", \"72\u00f76=\"],\n  [\"66\u00f77=\", \"90\u00f75=\"],\n  [\"54\u00f73=\", \"29\u00f79=\"],\n  [\"11\u00f74=\", \"83\u00f79=\"],\n  [\"43\u00f73=\", \"12\u00f72=\"],\n  [\"72\u00f77=\", \"35\u00f76=\"],\n  [\"75\u00f77=\", \"81\u00f74=\"],\n  [\"68\u00f74=\", \"26\u00f72=\"],\n  [\"34\u00f78=\", \"55\u00f78=\"],\n  [\"18\u00f78=\", \"91\u00f75=\"],\n  [\"44\u00f72=\", \"19\u00f75=\"],\n  [\"35\u00f73=\", \"60\u00f78=\"],\n  [\"67\u00f72=\", \"88\u00f73=\"],\n  [\"18\u00f73=\", \"55\u00f78=\"],\n  [\"99\u00f72=\", \"11\u00f77=\"],\n  [\"88\u00f79=\", \"52\u00f75=\"],\n  [\"91\u00f74=\", \"26\u00f79=\"],\n  [\"62\u00f77=\", \"68\u00f76=\"],\n  [\"32\u00f77=\", \"84\u00f78=\"],\n  [\"11\u00f72=\", \"97\u00f78=\"],\n  [\"59\u00f72=\", \"32\u00f73=\"],\n  [\"25\u00f77=\", \"53\u00f79=\"],\n  [\"17\u00f74=\", \"86\u00f72=\"],\n  [\"11\u00f76=\", \"24\u00f72=\"],\n  [\"24\u00f76=\", \"46\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Only the first occurrence should exist (all search strings are\n  // unique in the source document), but guard against duplicates by\n  // replacing every match found, just in case.\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the division problems in the table\n# with their updated values, as described by the diff. Each lookup\n# string is unique in the document, so a plain Find/Replace for each\n# pair is safe and unambiguous (even though some replacement values\n# repeat, e.g. \"55\u00f78=\" appears twice).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-08-05 Tuesday\", \"2025-08-06 Wednesday\"),\n  @(\"40\u00f74=\", \"72\u00f76=\"),\n  @(\"66\u00f77=\", \"90\u00f75=\"),\n  @(\"54\u00f73=\", \"29\u00f79=\"),\n  @(\"11\u00f74=\", \"83\u00f79=\"),\n  @(\"43\u00f73=\", \"12\u00f72=\"),\n  @(\"72\u00f77=\", \"35\u00f76=\"),\n  @(\"75\u00f77=\", \"81\u00f74=\"),\n  @(\"68\u00f74=\", \"26\u00f72=\"),\n  @(\"34\u00f78=\", \"55\u00f78=\"),\n  @(\"18\u00f78=\", \"91\u00f75=\"),\n  @(\"44\u00f72=\", \"19\u00f75=\"),\n  @(\"35\u00f73=\", \"60\u00f78=\"),\n  @(\"67\u00f72=\", \"88\u00f73=\"),\n  @(\"18\u00f73=\", \"55\u00f78=\"),\n  @(\"99\u00f72=\", \"11\u00f77=\"),\n  @(\"88\u00f79=\", \"52\u00f75=\"),\n  @(\"91\u00f74=\", \"26\u00f79=\"),\n  @(\"62\u00f77=\", \"68\u00f76=\"),\n  @(\"32\u00f77=\", \"84\u00f78=\"),\n  @(\"11\u00f72=\", \"97\u00f78=\"),\n  @(\"59\u00f72=\", \"32\u00f73=\"),\n  @(\"25\u00f77=\", \"53\u00f79=\"),\n  @(\"17\u00f74=\", \"86\u00f72=\"),\n  @(\"11\u00f76=\", \"24\u00f72=\"),\n  @(\"24\u00f76=\", \"46\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
